$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in row 5 with a new review row (mirrors the existing rows 2-4) ---
$ws.Range("A5").Value = "com.singleton.helix"
$ws.Range("B5").Value = "helix"
$ws.Range("C5").Value = "gregneri12@gmail.com"
$ws.Range("D5").Value = "halachme@gmail.com"
$ws.Range("E5").Value = "27/5/2019 15:56"
$ws.Range("F5").Value = "Welcome to the space jungle of helix jump! Haha great game."

# --- Hyperlink the two e-mail cells, same as the other review rows ---
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:gregneri12@gmail.com", "", "", "gregneri12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:halachme@gmail.com", "", "", "halachme@gmail.com")

# --- Restore the look of the populated cells to match the other rows ---
# (Hyperlinks.Add above stamps its own default hyperlink formatting; copy the
# formatting actually used by row 2's e-mail cells back over it.)
$ws.Range("C2").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D5").PasteSpecial(-4122)

# --- Row 5 was an empty placeholder row with a taller height; once it holds
# real data it should match the other data rows' (default) height. ---
$ws.Rows.Item(5).AutoFit()

# --- The placeholder "*unknown*" cell style that only the now-overwritten
# empty hyperlink-look cells used is no longer referenced; drop it. ---
$wb.Styles.Item("*unknown*").Delete()

$ws.Application.CutCopyMode = $false

# --- Match the author's final cursor position. ---
$ws.Range("C5").Select() | Out-Null
